$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update URL, Version, Date, Publisher values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-snapshot-age-in-years"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Elements": clear the Constraint(s) value for the root "Extension" row (row 2, column AI) ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""

# The "Fixed Value" for Extension.url (row 5) mirrors the StructureDefinition URL, so update it too
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-snapshot-age-in-years"
